$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 1.626987699542094
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 5.582307763322248

$ws.Range("B3").Value = 0.6545652718822623
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 3.223369029078222
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.038307959104277

$ws.Range("B4").Value = 0.2881169905109251
$ws.Range("C4").Value = 0.04103571897497393
$ws.Range("D4").Value = 0.7210945179870265
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 1.583633186074624
